$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: Typo in default Specific_ions config file
# Rows 14/20/31 are POS-mode / [M+H]+ fragments that were mistakenly
# saved with a trailing "-" (negative-ion) formula instead of "+".
$ws.Range("D31").Value = "C3H9NO6P+"
$ws.Range("D14").Value = "C5H15NO4P+"
$ws.Range("D20").Value = "C2H9O4NP+"

# Row 34's FORMULA cell held the literal arithmetic "NH3+H2O" instead of
# the resulting molecular formula "NH5O".
$ws.Range("D34").Value = "NH5O"

# Reflect the author's last selection (cell D34) when the file was saved.
$ws.Range("D34").Select()
